$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 3886
$ws.Range("F4").Value = 2307
$ws.Range("F9").Value = 110
$ws.Range("F10").Value = 111
$ws.Range("F11").Value = 1440
$ws.Range("F13").Value = 2564

# Sheet "演出" (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 35

# Sheet "全部类型" (All types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 3886
$ws.Range("F4").Value = 2307
$ws.Range("F8").Value = 35
$ws.Range("F10").Value = 110
$ws.Range("F11").Value = 111
$ws.Range("F14").Value = 1440
$ws.Range("F16").Value = 2564
